# updates to CIViC and MitelmanDB
#
# - Bump the Mitelman Database "source_version" value (biomarkers!E3) from
#   v20241015 to v20250115.
# - The workbook was left with the "biomarkers" tab active/selected (cursor
#   parked at C8), rather than the "compounds" tab.

$wb = $excel.ActiveWorkbook

$wsBiomarkers = $wb.Worksheets.Item("biomarkers")

# Mitelman Database source_version bump.
$wsBiomarkers.Range("E3").Value = "v20250115"

# Make "biomarkers" the active sheet and move the cursor to C8, matching the
# saved workbook/UI state.
$wsBiomarkers.Activate()
$wsBiomarkers.Range("C8").Select()
